$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 141.58333
$ws.Range("I9").Value = 155.44444
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 155.44444
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = 13.55556000000001
$ws.Range("N9").Value = -438

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 270.69232
$ws.Range("I33").Value = 166.66667
$ws.Range("J33").Value = 504.75
$ws.Range("K33").Value = 166.66667
$ws.Range("L33").Value = 504.75
$ws.Range("M33").Value = 62.33332999999999
$ws.Range("N33").Value = -962.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 645.82355
$ws.Range("I41").Value = 416.5
$ws.Range("J41").Value = 770.9091
$ws.Range("K41").Value = 416.5
$ws.Range("L41").Value = 770.9091
$ws.Range("M41").Value = 23.5
$ws.Range("N41").Value = -1650.9091

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 250
$ws.Range("I55").Value = 200
$ws.Range("J55").Value = 275
$ws.Range("K55").Value = 200
$ws.Range("L55").Value = 275
$ws.Range("M55").Value = 14
$ws.Range("N55").Value = -703

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2525.32
$ws.Range("I138").Value = 888
$ws.Range("J138").Value = 3071.0933
$ws.Range("K138").Value = 2664
$ws.Range("L138").Value = 9213.2799
$ws.Range("M138").Value = 2476
$ws.Range("N138").Value = -19493.2799

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5498.245
$ws.Range("I141").Value = 5717.476
$ws.Range("J141").Value = 4182.857
$ws.Range("K141").Value = 17152.428
$ws.Range("L141").Value = 12548.571
$ws.Range("M141").Value = -11972.428
$ws.Range("N141").Value = -22908.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3996.516
$ws.Range("I32").Value = 3516.1667
$ws.Range("J32").Value = 18407
$ws.Range("K32").Value = 3516.1667
$ws.Range("L32").Value = 18407
$ws.Range("M32").Value = -3229.1667
$ws.Range("N32").Value = -18981

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3134.9524
$ws.Range("I74").Value = 3119.5881
$ws.Range("J74").Value = 3200.25
$ws.Range("K74").Value = 3119.5881
$ws.Range("L74").Value = 3200.25
$ws.Range("M74").Value = -2245.5881
$ws.Range("N74").Value = -4948.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3134.9524
$ws.Range("I77").Value = 3119.5881
$ws.Range("J77").Value = 3200.25
$ws.Range("K77").Value = 15597.9405
$ws.Range("L77").Value = 16001.25
$ws.Range("M77").Value = -11229.9405
$ws.Range("N77").Value = -24737.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1730.8806
$ws.Range("I132").Value = 1173.6792
$ws.Range("J132").Value = 3840.2856
$ws.Range("K132").Value = 3521.0376
$ws.Range("L132").Value = 11520.8568
$ws.Range("M132").Value = -991.0376000000001
$ws.Range("N132").Value = -16580.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 43230.715
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 43230.715
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 43230.715
$ws.Range("N139").Value = -53510.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7694656.5
$ws.Range("I31").Value = 1226.3611
$ws.Range("J31").Value = 17245122
$ws.Range("K31").Value = 1226.3611
$ws.Range("L31").Value = 17245122
$ws.Range("M31").Value = -931.3611000000001
$ws.Range("N31").Value = -17245712

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7694656.5
$ws.Range("I34").Value = 1226.3611
$ws.Range("J34").Value = 17245122
$ws.Range("K34").Value = 1226.3611
$ws.Range("L34").Value = 17245122
$ws.Range("M34").Value = -1024.3611
$ws.Range("N34").Value = -17245526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2995.2903
$ws.Range("I132").Value = 2437.6785
$ws.Range("J132").Value = 8199.666999999999
$ws.Range("K132").Value = 7313.0355
$ws.Range("L132").Value = 24599.001
$ws.Range("M132").Value = -4783.0355
$ws.Range("N132").Value = -29659.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 13217.5
$ws.Range("I11").Value = 17358.572
$ws.Range("J11").Value = 7420
$ws.Range("K11").Value = 52075.716
$ws.Range("L11").Value = 22260
$ws.Range("M11").Value = -51935.716
$ws.Range("N11").Value = -22540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 2354.5293
$ws.Range("I21").Value = 607.25
$ws.Range("J21").Value = 2892.1538
$ws.Range("K21").Value = 1821.75
$ws.Range("L21").Value = 8676.4614
$ws.Range("M21").Value = -1648.75
$ws.Range("N21").Value = -9022.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1682.909
$ws.Range("I70").Value = 930.2857
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 2790.8571
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -2475.8571
$ws.Range("N70").Value = -9630

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 1682.909
$ws.Range("I73").Value = 930.2857
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 2790.8571
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -1698.8571
$ws.Range("N73").Value = -11184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1200
$ws.Range("I75").Value = 100
$ws.Range("J75").Value = 1750
$ws.Range("K75").Value = 300
$ws.Range("L75").Value = 5250
$ws.Range("M75").Value = 698
$ws.Range("N75").Value = -7246

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 1200
$ws.Range("I78").Value = 100
$ws.Range("J78").Value = 1750
$ws.Range("K78").Value = 900
$ws.Range("L78").Value = 15750
$ws.Range("M78").Value = 4092
$ws.Range("N78").Value = -25734

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 24209.715
$ws.Range("I107").Value = 449.41177
$ws.Range("J107").Value = 40366.72
$ws.Range("K107").Value = 1348.23531
$ws.Range("L107").Value = 121100.16
$ws.Range("M107").Value = 571.76469
$ws.Range("N107").Value = -124940.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11905642
$ws.Range("I131").Value = 71429130
$ws.Range("J131").Value = 945.6
$ws.Range("K131").Value = 214287390
$ws.Range("L131").Value = 2836.8
$ws.Range("M131").Value = -214282350
$ws.Range("N131").Value = -12916.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3888.9487
$ws.Range("I134").Value = 2965.95
$ws.Range("J134").Value = 4860.5264
$ws.Range("K134").Value = 8897.849999999999
$ws.Range("L134").Value = 14581.5792
$ws.Range("M134").Value = -3827.849999999999
$ws.Range("N134").Value = -24721.5792

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10872053
$ws.Range("I80").Value = 50002084
$ws.Range("J80").Value = 2600
$ws.Range("K80").Value = 50002084
$ws.Range("L80").Value = 2600
$ws.Range("M80").Value = -50001086
$ws.Range("N80").Value = -4596

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 10872053
$ws.Range("I83").Value = 50002084
$ws.Range("J83").Value = 2600
$ws.Range("K83").Value = 250010420
$ws.Range("L83").Value = 13000
$ws.Range("M83").Value = -250005428
$ws.Range("N83").Value = -22984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1588.5714
$ws.Range("I102").Value = 1074.2
$ws.Range("J102").Value = 2274.4
$ws.Range("K102").Value = 1074.2
$ws.Range("L102").Value = 2274.4
$ws.Range("M102").Value = 547.8
$ws.Range("N102").Value = -5518.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4156.433
$ws.Range("I122").Value = 2436.0527
$ws.Range("J122").Value = 7128
$ws.Range("K122").Value = 7308.158100000001
$ws.Range("L122").Value = 21384
$ws.Range("M122").Value = -4858.158100000001
$ws.Range("N122").Value = -26284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11790.485
$ws.Range("I132").Value = 15899.823
$ws.Range("J132").Value = 7909.4443
$ws.Range("K132").Value = 47699.469
$ws.Range("L132").Value = 23728.3329
$ws.Range("M132").Value = -45169.469
$ws.Range("N132").Value = -28788.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 40000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 40000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 40000
$ws.Range("N48").Value = -41138

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1930.3396
$ws.Range("I136").Value = 615.3684
$ws.Range("J136").Value = 5261.6
$ws.Range("K136").Value = 1846.1052
$ws.Range("L136").Value = 15784.8
$ws.Range("M136").Value = 703.8948
$ws.Range("N136").Value = -20884.8
